# Commit: "Add copyrights for figures in Chapter 2"
#
# The thesis correction-list sheet gets answers filled in for the
# "Chapter 2 / figure permissions" comments: a new answer in C4, a small
# wording fix in A5, and a brand new Q/A pair (B17/C17) about obtaining
# reproduction permission for Figure 2.5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C4: answer for "2) Chapter 2: Make sure you have permission..."
$ws.Range("C4").Value = "Permissions have been obtain. The copyrights of figures are stated in Figure 2.1, 2.4, 2.5, 2.6, 2.7, 2.8, 2.9 following the rules of publishers. "

# A5: small wording/punctuation fix ("possible" -> "possible.")
$ws.Range("A5").Value = "3) Chapter 3: Add a section to describe the experimental methodology in detail, including figures of the setup if possible. Include a method section to show how you obtained the data and a picture of the experiment."

# B17/C17: new note + answer about permission for figure reproduction (replaces
# the old, now-duplicated "natural log" note which was already covered in row 21)
$ws.Range("B17").Value = "obtain permission for any image reproductions such as fig 2.5 This is usually just a matter of applying to the publisher. See https://www.acm.org/publications/policies/copyright-policy "
$ws.Range("C17").Value = "Permissions have been obtain. The copyrights of figures are stated in Figure 2.1, 2.4, 2.5, 2.6, 2.7, 2.8, 2.9 following the rules of publishers (same in comment 2)."

# Row 17 now holds much longer text, so grow the row to fit it (matches the
# ht="74" seen after the edit).
$ws.Rows("17:17").RowHeight = 74

# The author's cursor ended up one row lower (C4 -> C5) after filling these in.
$ws.Range("C5").Select()
